# feat: stabilize member-consult-measure workflow with korean columns
#
# raw_stock: append 4 new "자동 원단 소요 처리" (auto fabric usage) OUT rows
# for stock F001, and roll the F001 monthly usage / balance aggregates
# forward to reflect the new consumption.

$wb = $excel.ActiveWorkbook

$wsStock = $wb.Worksheets.Item("raw_stock")

# Column A on the existing data rows carries a date number format (the
# "date" column) - reuse it for the newly appended rows instead of
# hard-coding a format string.
$dateFmt = $wsStock.Range("A2").NumberFormat()

$newRows = @(
    @{ Row = 8;  Date = 46002; StockId = "F001"; StockName = "이태리 순모 네이비"; Type = "OUT"; Qty = 2.7;  Unit = "m"; OrderId = "2025-3811-01";  Note = "자동 원단 소요 처리"; QtySigned = -2.7; Month = "2025-12" }
    @{ Row = 9;  Date = 46020; StockId = "F001"; StockName = "이태리 순모 네이비"; Type = "OUT"; Qty = 2.7;  Unit = "m"; OrderId = "2025-10000-01"; Note = "자동 원단 소요 처리"; QtySigned = -2.7; Month = "2025-12" }
    @{ Row = 10; Date = 46020; StockId = "F001"; StockName = "이태리 순모 네이비"; Type = "OUT"; Qty = 2.7;  Unit = "m"; OrderId = "2025-10000-01"; Note = "자동 원단 소요 처리"; QtySigned = -2.7; Month = "2025-12" }
    @{ Row = 11; Date = 46020; StockId = "F001"; StockName = "이태리 순모 네이비"; Type = "OUT"; Qty = 2.7;  Unit = "m"; OrderId = "2025-10000-02"; Note = "자동 원단 소요 처리"; QtySigned = -2.7; Month = "2025-12" }
)

foreach ($row in $newRows) {
    $r = $row.Row
    $wsStock.Range("A$r").Value = $row.Date
    $wsStock.Range("B$r").Value = $row.StockId
    $wsStock.Range("C$r").Value = $row.StockName
    $wsStock.Range("D$r").Value = $row.Type
    $wsStock.Range("E$r").Value = $row.Qty
    $wsStock.Range("F$r").Value = $row.Unit
    $wsStock.Range("G$r").Value = $row.OrderId
    $wsStock.Range("H$r").Value = $row.Note
    $wsStock.Range("I$r").Value = $row.QtySigned
    $wsStock.Range("J$r").Value = $row.Month

    $wsStock.Range("A$r").NumberFormat = $dateFmt
}

# usage: roll the F001 / 2025-12 quantity_signed total forward by the
# 4 x -2.7 consumed above (30 -> 19.2).
$wsUsage = $wb.Worksheets.Item("usage")
$wsUsage.Range("C5").Value = 19.2

# balance: roll the F001 running balance forward by the same amount
# (27.7 -> 16.9).
$wsBalance = $wb.Worksheets.Item("balance")
$wsBalance.Range("B3").Value = 16.9
